$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24; this shifts existing rows 24:135 down to 25:136
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record (same shape/style as its
# neighbours: Mango / Terminal Hortofrutícola Agro Chillán / Ñuble)
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = (Get-Date -Year 2023 -Month 5 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108002
$ws.Range("J24").Value = "Mango"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 8000
$ws.Range("O24").Value = 8000
$ws.Range("P24").Value = 8000
$ws.Range("Q24").Value = "`$/bandeja 4 kilos"
$ws.Range("R24").Value = "Perú"
$ws.Range("S24").Value = 2000
$ws.Range("T24").Value = 4
